$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$oldText = $wsHoja1.Range("A1").Value2
$newText = $oldText `
    -replace [regex]::Escape("1000 Bs = 1.57 = 5952.13 pesos"), "1000 Bs = 1.5 = 5613.83 pesos" `
    -replace [regex]::Escape("5952.13 pesos = 1.58 = 924.75 Bs"), "5613.83 pesos = 1.49 = 925.69 Bs"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate figures in N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 665.499
$wsTasas.Range("O10").Value = 3736
$wsTasas.Range("N12").Value = 3760
$wsTasas.Range("O12").Value = 620
